$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the updated Price/Volume cells to remain plain text (matches the source data,
# which is non-numeric: thousand-dot formatted prices and "  +x.xx%  " strings) so
# Excel does not silently reinterpret them as numbers and drop formatting/precision.
$textCells = @("D2", "D3", "D5", "D6", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D27", "D28", "D30", "D31", "D32", "D35", "D36", "D37", "D39", "D42", "D43", "D46", "D48", "E2", "E3", "E4", "E5", "E6", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16", "E17", "E18", "E19", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E27", "E28", "E29", "E30", "E31", "E32", "E33", "E34", "E35", "E36", "E37", "E38", "E39", "E40", "E41", "E42", "E43", "E44", "E45", "E46", "E48", "E49", "E50", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.998.21'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '3.129.47'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '580.77'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").Value = '173.61'
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '37.56'
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("D14").Value = '0.122'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").Value = '66.972.06'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '7.14'
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '3.127.82'
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").Value = '16.41'
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("D19").Value = '489.87'
$ws.Range("E19").Value = '  +1.38%  '
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").Value = '7.91'
$ws.Range("E21").Value = '  +5.34%  '
$ws.Range("D22").Value = '84.21'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '13.24'
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("E24").Value = '  -1.62%  '
$ws.Range("D25").Value = '10.35'
$ws.Range("E25").Value = '  +3.76%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = '7.95'
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").Value = '2.36'
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = '28.73'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").Value = '0.114'
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").Value = '0.0₃0950'
$ws.Range("E32").Value = '  -6.03%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = '0.978'
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("D36").Value = '46.90'
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").Value = '50.12'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("D39").Value = '0.312'
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("E40").Value = '  +1.76%  '
$ws.Range("E41").Value = '  -1.00%  '
$ws.Range("D42").Value = '386.80'
$ws.Range("E42").Value = '  +2.30%  '
$ws.Range("D43").Value = '2.825.93'
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("E44").Value = '  -6.85%  '
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("D46").Value = '135.99'
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D48").Value = '25.19'
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("E49").Value = '  +0.61%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("E51").Value = '  -0.25%  '

Write-Output "Updated 28 price cells and 48 volume cells in cryptos sheet."
